$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 21-34 with new values per diff ---
$ws.Range("D21").Value = 44567
$ws.Range("K21").Value = "Modesto"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 18000
$ws.Range("O21").Value = 18000
$ws.Range("P21").Value = 18000
$ws.Range("Q21").Value = "$/caja 18 kilos"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 1000
$ws.Range("T21").Value = 18

$ws.Range("D22").Value = 44567
$ws.Range("K22").Value = "Modesto"
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 250
$ws.Range("N22").Value = 13000
$ws.Range("O22").Value = 13000
$ws.Range("P22").Value = 13000
$ws.Range("Q22").Value = "$/caja 15 kilos"
$ws.Range("R22").Value = "Región de O'Higgins"
$ws.Range("S22").Value = 867
$ws.Range("T22").Value = 15

$ws.Range("D23").Value = 44176
$ws.Range("K23").Value = "Castle Brite"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 17000
$ws.Range("O23").Value = 17000
$ws.Range("P23").Value = 17000
$ws.Range("Q23").Value = "$/caja 18 kilos granel"
$ws.Range("R23").Value = "Provincia de Limarí"
$ws.Range("S23").Value = 944
$ws.Range("T23").Value = 18

$ws.Range("D24").Value = 44168
$ws.Range("K24").Value = "Castle Brite"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 250
$ws.Range("N24").Value = 10000
$ws.Range("O24").Value = 10000
$ws.Range("P24").Value = 10000
$ws.Range("Q24").Value = "$/caja 10 kilos"
$ws.Range("R24").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S24").Value = 1000
$ws.Range("T24").Value = 10

$ws.Range("D25").Value = 44168
$ws.Range("K25").Value = "Castle Brite"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 100
$ws.Range("N25").Value = 17000
$ws.Range("O25").Value = 17000
$ws.Range("P25").Value = 17000
$ws.Range("Q25").Value = "$/caja 18 kilos"
$ws.Range("R25").Value = "Provincia de Limarí"
$ws.Range("S25").Value = 944
$ws.Range("T25").Value = 18

$ws.Range("D26").Value = 44553
$ws.Range("K26").Value = "Castle Brite"
$ws.Range("L26").Value = "Especial"
$ws.Range("M26").Value = 100
$ws.Range("N26").Value = 13000
$ws.Range("O26").Value = 13000
$ws.Range("P26").Value = 13000
$ws.Range("Q26").Value = "$/caja 10 kilos"
$ws.Range("R26").Value = "Región de O'Higgins"
$ws.Range("S26").Value = 1300
$ws.Range("T26").Value = 10

$ws.Range("D27").Value = 44553
$ws.Range("K27").Value = "Castle Brite"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 120
$ws.Range("N27").Value = 11000
$ws.Range("O27").Value = 11000
$ws.Range("P27").Value = 11000
$ws.Range("Q27").Value = "$/caja 10 kilos"
$ws.Range("R27").Value = "Región de O'Higgins"
$ws.Range("S27").Value = 1100
$ws.Range("T27").Value = 10

$ws.Range("D28").Value = 44553
$ws.Range("K28").Value = "Castle Brite"
$ws.Range("L28").Value = "Segunda"
$ws.Range("M28").Value = 150
$ws.Range("N28").Value = 14000
$ws.Range("O28").Value = 14000
$ws.Range("P28").Value = 14000
$ws.Range("Q28").Value = "$/caja 15 kilos"
$ws.Range("R28").Value = "Región de O'Higgins"
$ws.Range("S28").Value = 933
$ws.Range("T28").Value = 15

$ws.Range("D29").Value = 44167
$ws.Range("K29").Value = "Castle Brite"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 300
$ws.Range("N29").Value = 15000
$ws.Range("O29").Value = 15000
$ws.Range("P29").Value = 15000
$ws.Range("Q29").Value = "$/caja 16 kilos granel"
$ws.Range("R29").Value = "Provincia de Limarí"
$ws.Range("S29").Value = 938
$ws.Range("T29").Value = 16

$ws.Range("D30").Value = 44161
$ws.Range("K30").Value = "Castle Brite"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 150
$ws.Range("N30").Value = 20000
$ws.Range("O30").Value = 20000
$ws.Range("P30").Value = 20000
$ws.Range("Q30").Value = "$/caja 18 kilos granel"
$ws.Range("R30").Value = "Provincia de Limarí"
$ws.Range("S30").Value = 1111
$ws.Range("T30").Value = 18

$ws.Range("D31").Value = 44160
$ws.Range("K31").Value = "Dina"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 200
$ws.Range("N31").Value = 20000
$ws.Range("O31").Value = 20000
$ws.Range("P31").Value = 20000
$ws.Range("Q31").Value = "$/caja 15 kilos"
$ws.Range("R31").Value = "Provincia de Limarí"
$ws.Range("S31").Value = 1333
$ws.Range("T31").Value = 15

$ws.Range("D32").Value = 44543
$ws.Range("K32").Value = "Castle Brite"
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 100
$ws.Range("N32").Value = 18000
$ws.Range("O32").Value = 18000
$ws.Range("P32").Value = 18000
$ws.Range("Q32").Value = "$/caja 15 kilos"
$ws.Range("R32").Value = "Región de O'Higgins"
$ws.Range("S32").Value = 1200
$ws.Range("T32").Value = 15

$ws.Range("D33").Value = 44543
$ws.Range("K33").Value = "Castle Brite"
$ws.Range("L33").Value = "Segunda"
$ws.Range("M33").Value = 50
$ws.Range("N33").Value = 15000
$ws.Range("O33").Value = 15000
$ws.Range("P33").Value = 15000
$ws.Range("Q33").Value = "$/caja 15 kilos"
$ws.Range("R33").Value = "Región de O'Higgins"
$ws.Range("S33").Value = 1000
$ws.Range("T33").Value = 15

$ws.Range("D34").Value = 44187
$ws.Range("K34").Value = "Dina"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 120
$ws.Range("N34").Value = 16000
$ws.Range("O34").Value = 16000
$ws.Range("P34").Value = 16000
$ws.Range("Q34").Value = "$/caja 18 kilos"
$ws.Range("R34").Value = "Provincia de Limarí"
$ws.Range("S34").Value = 889
$ws.Range("T34").Value = 18

# --- Insert two new rows at position 35, pushing old row 35 down to row 37 ---
$ws.Rows.Item(35).Insert()
$ws.Rows.Item(36).Insert()

# Row 35
$ws.Range("A35").Value = 5
$ws.Range("B35").Value = "Macroferia Regional de Talca"
$ws.Range("C35").Value = "Maule"
$ws.Range("D35").Value = 44529
$ws.Range("E35").Value = 7
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100103
$ws.Range("H35").Value = "Frutos de hueso (carozo)"
$ws.Range("I35").Value = 100103003
$ws.Range("J35").Value = "Damasco"
$ws.Range("K35").Value = "Castle Brite"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 100
$ws.Range("N35").Value = 20000
$ws.Range("O35").Value = 20000
$ws.Range("P35").Value = 20000
$ws.Range("Q35").Value = "$/caja 15 kilos"
$ws.Range("R35").Value = "Región de O'Higgins"
$ws.Range("S35").Value = 1333
$ws.Range("T35").Value = 15

# Row 36
$ws.Range("A36").Value = 5
$ws.Range("B36").Value = "Macroferia Regional de Talca"
$ws.Range("C36").Value = "Maule"
$ws.Range("D36").Value = 44540
$ws.Range("E36").Value = 7
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100103
$ws.Range("H36").Value = "Frutos de hueso (carozo)"
$ws.Range("I36").Value = 100103003
$ws.Range("J36").Value = "Damasco"
$ws.Range("K36").Value = "Castle Brite"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 120
$ws.Range("N36").Value = 20000
$ws.Range("O36").Value = 20000
$ws.Range("P36").Value = 20000
$ws.Range("Q36").Value = "$/caja 16 kilos"
$ws.Range("R36").Value = "Región de O'Higgins"
$ws.Range("S36").Value = 1250
$ws.Range("T36").Value = 16

"Done"